$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.327.09'
$ws.Range("E2").Value = '  +0.66%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.881.59'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  +0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '604.41'
$ws.Range("E5").Value = '  +0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '168.30'
$ws.Range("E6").Value = '  +3.62%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.859.60'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +0.36%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.533'
$ws.Range("E9").Value = '  +0.34%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.168'
$ws.Range("E10").Value = '  -0.35%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.37'
$ws.Range("E11").Value = '  +1.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.465'
$ws.Range("E12").Value = '  +1.39%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000253'
$ws.Range("E13").Value = '  +4.05%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.97'
$ws.Range("E14").Value = '  +2.70%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.543.91'
$ws.Range("E15").Value = '  +0.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.898.10'
$ws.Range("E16").Value = '  +0.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '69.537.42'
$ws.Range("E17").Value = '  +0.72%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.65'
$ws.Range("E18").Value = '  +9.17%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.57'
$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("E20").Value = '  -1.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.10'
$ws.Range("E21").Value = '  -2.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.41'
$ws.Range("E22").Value = '  +0.55%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.740'
$ws.Range("E23").Value = '  +2.94%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000167'
$ws.Range("E24").Value = '  +2.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.28'
$ws.Range("E25").Value = '  +1.48%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  +1.40%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.26'
$ws.Range("E27").Value = '  +1.24%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.10'
$ws.Range("E28").Value = '  +1.32%  '

$ws.Range("E29").Value = '  +0.05%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.97'
$ws.Range("E30").Value = '  +0.61%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.040.03'
$ws.Range("E31").Value = '  +0.60%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.41'
$ws.Range("E32").Value = '  +1.80%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.77'
$ws.Range("E33").Value = '  -1.91%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '31.84'
$ws.Range("E34").Value = '  -1.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.851.09'
$ws.Range("E35").Value = '  +0.99%  '

$ws.Range("E36").Value = '  -0.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.07'
$ws.Range("E37").Value = '  +3.05%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.03'
$ws.Range("E38").Value = '  -0.09%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.141'
$ws.Range("E39").Value = '  +1.00%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.27'
$ws.Range("E40").Value = '  +10.06%  '

$ws.Range("E41").Value = '  +0.46%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.327'
$ws.Range("E42").Value = '  +2.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.11'
$ws.Range("E43").Value = '  +6.22%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '434.61'
$ws.Range("E44").Value = '  -0.62%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.00'
$ws.Range("E45").Value = '  -1.01%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.62'
$ws.Range("E46").Value = '  +2.54%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0367'
$ws.Range("E48").Value = '  +2.44%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '143.04'
$ws.Range("E49").Value = '  -0.28%  '

$ws.Range("B50").Value = 'Arweave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '40.18'
$ws.Range("E50").Value = '  +3.49%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000269'
$ws.Range("E51").Value = '  +18.28%  '
